# #15: Container comps are not react-agnostic
#
# - "TextBox 7" (the Container "glue" blurb): reword + reflow/resize.
# - "TextBox 8" (the Component blurb): drop the stray "React" + resize.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- TextBox 7 ("Glue" blurb, under the Container box) ---------------------
$sh3 = $s.Shapes.Item(3)

# Reposition / resize (a:off/a:ext in EMU -> points, 12700 EMU per point).
$sh3.Left  = 417.5972450944882
$sh3.Width = 214.6711053622047

# Reword while preserving the existing soft line-break between the two runs:
# first run becomes "Glue & Logic", second run becomes the former title text.
$tr3 = $sh3.TextFrame.TextRange
$full3 = $tr3.Text
$brk3 = $full3.IndexOf([char]11)

$run3a = $tr3.Characters(1, $brk3)
$run3a.Text = "Glue & Logic"

$full3b = $tr3.Text
$brk3b = $full3b.IndexOf([char]11)
$run3b = $tr3.Characters($brk3b + 2, $full3b.Length - $brk3b - 1)
$run3b.Text = "Abstract over Redux state handling"

# --- TextBox 8 (blurb under the Component box) ------------------------------
$sh8 = $s.Shapes.Item(8)

# Reposition / resize (a:off/a:ext in EMU -> points, 12700 EMU per point).
$sh8.Left  = 768.0430608661418
$sh8.Width = 118.68283464566929

# Drop the trailing "React" from the first line; second line is unchanged.
$tr8 = $sh8.TextFrame.TextRange
$full8 = $tr8.Text
$brk8 = $full8.IndexOf([char]11)

$run8a = $tr8.Characters(1, $brk8)
$run8a.Text = "Presentation-Only"
